$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    17  = "G0S2"
    21  = "ZNF503"
    41  = "KIF26A"
    52  = "VAT1L"
    54  = "CENPV"
    56  = "C1QL1"
    65  = "CTXN1"
    71  = "INHBB"
    85  = "OLIG1"
    100 = "FGFBP2"
    101 = "SLC10A4"
    105 = "C5orf49"
    115 = "NSG2"
    119 = "ID4"
    140 = "GABBR2"
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
